# Apply updated cryptocurrency data to the worksheet, matching the source diff.
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range('D2').NumberFormat = '@'
$ws.Range('D2').Value = '45.049.08'
$ws.Range('E2').Value = '  +2.79%  '
$ws.Range('D3').NumberFormat = '@'
$ws.Range('D3').Value = '2.357.52'
$ws.Range('E3').Value = '  +0.91%  '
$ws.Range('E4').Value = '  -0.16%  '
$ws.Range('D5').NumberFormat = '@'
$ws.Range('D5').Value = '311.20'
$ws.Range('E5').Value = '  -0.24%  '
$ws.Range('D6').NumberFormat = '@'
$ws.Range('D6').Value = '107.47'
$ws.Range('E6').Value = '  -0.77%  '
$ws.Range('E7').Value = '  -0.34%  '
$ws.Range('E8').Value = '  -0.22%  '
$ws.Range('E9').Value = '  -2.03%  '
$ws.Range('D10').NumberFormat = '@'
$ws.Range('D10').Value = '40.68'
$ws.Range('E10').Value = '  -1.32%  '
$ws.Range('D11').NumberFormat = '@'
$ws.Range('D11').Value = '0.0915'
$ws.Range('E11').Value = '  -0.22%  '
$ws.Range('D12').NumberFormat = '@'
$ws.Range('D12').Value = '8.41'
$ws.Range('E12').Value = '  -1.52%  '
$ws.Range('E13').Value = '  +1.18%  '
$ws.Range('D14').NumberFormat = '@'
$ws.Range('D14').Value = '0.974'
$ws.Range('E14').Value = '  -3.16%  '
$ws.Range('D15').NumberFormat = '@'
$ws.Range('D15').Value = '2.714.95'
$ws.Range('E15').Value = '  +0.73%  '
$ws.Range('D16').NumberFormat = '@'
$ws.Range('D16').Value = '15.14'
$ws.Range('E16').Value = '  -2.00%  '
$ws.Range('D17').NumberFormat = '@'
$ws.Range('D17').Value = '2.363.32'
$ws.Range('E17').Value = '  +1.27%  '
$ws.Range('D18').NumberFormat = '@'
$ws.Range('D18').Value = '44.986.78'
$ws.Range('E18').Value = '  +2.70%  '
$ws.Range('D19').NumberFormat = '@'
$ws.Range('D19').Value = '14.50'
$ws.Range('E19').Value = '  +11.10%  '
$ws.Range('B20').Value = 'ShibaInu'
$ws.Range('C20').Value = 'https://coinranking.com/coin/xz24e0BjL+shibainu-shib'
$ws.Range('D20').NumberFormat = '@'
$ws.Range('D20').Value = '0.0000106'
$ws.Range('E20').Value = '  -0.81%  '
$ws.Range('B21').Value = 'Uniswap'
$ws.Range('C21').Value = 'https://coinranking.com/coin/_H5FVG9iW+uniswap-uni'
$ws.Range('D21').NumberFormat = '@'
$ws.Range('D21').Value = '7.18'
$ws.Range('E21').Value = '  -4.81%  '
$ws.Range('D22').NumberFormat = '@'
$ws.Range('D22').Value = '72.85'
$ws.Range('E22').Value = '  -1.75%  '
$ws.Range('D23').NumberFormat = '@'
$ws.Range('D23').Value = '3.50'
$ws.Range('E23').Value = '  +1.10%  '
$ws.Range('D24').NumberFormat = '@'
$ws.Range('D24').Value = '258.03'
$ws.Range('E24').Value = '  -4.02%  '
$ws.Range('E25').Value = '  +1.42%  '
$ws.Range('E26').Value = '  -0.14%  '
$ws.Range('E27').Value = '  -0.52%  '
$ws.Range('D28').NumberFormat = '@'
$ws.Range('D28').Value = '7.15'
$ws.Range('E28').Value = '  -6.78%  '
$ws.Range('E29').Value = '  +1.11%  '
$ws.Range('D30').NumberFormat = '@'
$ws.Range('D30').Value = '0.0968'
$ws.Range('E30').Value = '  +9.16%  '
$ws.Range('D31').NumberFormat = '@'
$ws.Range('D31').Value = '22.23'
$ws.Range('E31').Value = '  -1.52%  '
$ws.Range('D32').NumberFormat = '@'
$ws.Range('D32').Value = '37.05'
$ws.Range('E32').Value = '  -5.19%  '
$ws.Range('D33').NumberFormat = '@'
$ws.Range('D33').Value = '167.59'
$ws.Range('E33').Value = '  -0.67%  '
$ws.Range('D34').NumberFormat = '@'
$ws.Range('D34').Value = '2.98'
$ws.Range('E34').Value = '  +5.39%  '
$ws.Range('E35').Value = '  -1.56%  '
$ws.Range('E36').Value = '  +1.06%  '
$ws.Range('D37').NumberFormat = '@'
$ws.Range('D37').Value = '4.66'
$ws.Range('E37').Value = '  -1.34%  '
$ws.Range('E38').Value = '  +4.36%  '
$ws.Range('D39').NumberFormat = '@'
$ws.Range('D39').Value = '2.89'
$ws.Range('E39').Value = '  +0.41%  '
$ws.Range('D40').NumberFormat = '@'
$ws.Range('D40').Value = '0.0350'
$ws.Range('E40').Value = '  -3.37%  '
$ws.Range('E41').Value = '  +2.38%  '
$ws.Range('D42').NumberFormat = '@'
$ws.Range('D42').Value = '99.94'
$ws.Range('E42').Value = '  -4.76%  '
$ws.Range('D43').NumberFormat = '@'
$ws.Range('D43').Value = '1.880.25'
$ws.Range('E43').Value = '  +12.40%  '
$ws.Range('D44').NumberFormat = '@'
$ws.Range('D44').Value = '69.05'
$ws.Range('E44').Value = '  -3.39%  '
$ws.Range('E45').Value = '  -4.37%  '
$ws.Range('B46').Value = 'FirstDigitalUSD'
$ws.Range('C46').Value = 'https://coinranking.com/coin/cpjRxjFYD+firstdigitalusd-fdusd'
$ws.Range('D46').NumberFormat = '@'
$ws.Range('D46').Value = '1.00'
$ws.Range('E46').Value = '  -0.40%  '
$ws.Range('B47').Value = 'Celestia'
$ws.Range('C47').Value = 'https://coinranking.com/coin/YQcD0lBl7+celestia-tia'
$ws.Range('D47').NumberFormat = '@'
$ws.Range('D47').Value = '12.78'
$ws.Range('E47').Value = '  -4.60%  '
$ws.Range('D48').NumberFormat = '@'
$ws.Range('D48').Value = '81.23'
$ws.Range('E48').Value = '  +5.89%  '
$ws.Range('D49').NumberFormat = '@'
$ws.Range('D49').Value = '5.61'
$ws.Range('E49').Value = '  +7.77%  '
$ws.Range('B50').Value = 'Aave'
$ws.Range('C50').Value = 'https://coinranking.com/coin/ixgUfzmLR+aave-aave'
$ws.Range('D50').NumberFormat = '@'
$ws.Range('D50').Value = '109.95'
$ws.Range('E50').Value = '  -3.75%  '
$ws.Range('B51').Value = 'FraxShare'
$ws.Range('C51').Value = 'https://coinranking.com/coin/3nNpuxHJ8+fraxshare-fxs'
$ws.Range('D51').NumberFormat = '@'
$ws.Range('D51').Value = '9.14'
$ws.Range('E51').Value = '  +2.31%  '
